$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nam "
$ws.Range("B1").Value = "jhjshdsa"
$ws.Range("C1").Value = "kjaskjdqaj"

$ws.Range("C1").Select()
